$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Denv 1 (row 4, NC_001477) major_lineage label corrected to drop the
# trailing space left over from the old reference set: "1IV " -> "1IV"
$ws.Range("C4").Value = "1IV"

# Leave selection where the editor ended up after making the change.
$ws.Range("G11").Select()
